$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("distancetable")

# Shift the whole distance table down by one row (insert a blank row at the top)
$ws.Rows.Item(1).Insert()

# Move selection / active cell to match the saved workbook state
$ws.Range("C19").Select()
